# ---------------------------------------------------------------------------
# Commit: "New .sql files and views, PowerBI Dashboard"
# Adds two new worksheets at the end of the workbook:
#   07_view_product_category_sales_   (sheetId 11) - product-category x month
#                                                     sales breakdown
#   08_view_sales_and_customers_mom   (sheetId 12) - duplicate export of the
#                                                     06_view_sales_and_customers_mom
#                                                     data (month-over-month revenue)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- New sheet: 07_view_product_category_sales_ ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws1.Name = "07_view_product_category_sales_"

# Column B (month_start) holds literal text dates like "2023-01-01". Force
# the column to Text format *before* writing so Excel does not coerce the
# strings into date serial numbers.
$ws1.Range("B1:B39").NumberFormat = "@"

    $ws1.Cells.Item(1, 1).Value = 'product_category'
    $ws1.Cells.Item(1, 2).Value = 'month_start'
    $ws1.Cells.Item(1, 3).Value = 'transactions_count'
    $ws1.Cells.Item(1, 4).Value = 'month_total_amount'
    $ws1.Cells.Item(1, 5).Value = 'month_total_quantity'
    $ws1.Cells.Item(1, 6).Value = 'pct_of_category_amount'
    $ws1.Cells.Item(1, 7).Value = 'pct_of_category_quantity'
    $ws1.Cells.Item(1, 8).Value = 'pct_of_month_amount'
    $ws1.Cells.Item(1, 9).Value = 'pct_of_month_quantity'
    $ws1.Cells.Item(2, 1).Value = 'Beauty'
    $ws1.Cells.Item(2, 2).Value = '2023-01-01'
    $ws1.Cells.Item(2, 3).Value = 25
    $ws1.Cells.Item(2, 4).Value = 12430
    $ws1.Cells.Item(2, 5).Value = 59
    $ws1.Cells.Item(2, 6).Value = 0.0866
    $ws1.Cells.Item(2, 7).Value = 0.0765
    $ws1.Cells.Item(2, 8).Value = 0.3506
    $ws1.Cells.Item(2, 9).Value = 0.3026
    $ws1.Cells.Item(3, 1).Value = 'Beauty'
    $ws1.Cells.Item(3, 2).Value = '2023-02-01'
    $ws1.Cells.Item(3, 3).Value = 26
    $ws1.Cells.Item(3, 4).Value = 14035
    $ws1.Cells.Item(3, 5).Value = 68
    $ws1.Cells.Item(3, 6).Value = 0.0978
    $ws1.Cells.Item(3, 7).Value = 0.0882
    $ws1.Cells.Item(3, 8).Value = 0.3185
    $ws1.Cells.Item(3, 9).Value = 0.3178
    $ws1.Cells.Item(4, 1).Value = 'Beauty'
    $ws1.Cells.Item(4, 2).Value = '2023-03-01'
    $ws1.Cells.Item(4, 3).Value = 21
    $ws1.Cells.Item(4, 4).Value = 10545
    $ws1.Cells.Item(4, 5).Value = 51
    $ws1.Cells.Item(4, 6).Value = 0.0735
    $ws1.Cells.Item(4, 7).Value = 0.06610000000000001
    $ws1.Cells.Item(4, 8).Value = 0.3637
    $ws1.Cells.Item(4, 9).Value = 0.2629
    $ws1.Cells.Item(5, 1).Value = 'Beauty'
    $ws1.Cells.Item(5, 2).Value = '2023-04-01'
    $ws1.Cells.Item(5, 3).Value = 29
    $ws1.Cells.Item(5, 4).Value = 11905
    $ws1.Cells.Item(5, 5).Value = 69
    $ws1.Cells.Item(5, 6).Value = 0.083
    $ws1.Cells.Item(5, 7).Value = 0.0895
    $ws1.Cells.Item(5, 8).Value = 0.3515
    $ws1.Cells.Item(5, 9).Value = 0.3224
    $ws1.Cells.Item(6, 1).Value = 'Beauty'
    $ws1.Cells.Item(6, 2).Value = '2023-05-01'
    $ws1.Cells.Item(6, 3).Value = 28
    $ws1.Cells.Item(6, 4).Value = 12450
    $ws1.Cells.Item(6, 5).Value = 65
    $ws1.Cells.Item(6, 6).Value = 0.0868
    $ws1.Cells.Item(6, 7).Value = 0.0843
    $ws1.Cells.Item(6, 8).Value = 0.2342
    $ws1.Cells.Item(6, 9).Value = 0.251
    $ws1.Cells.Item(7, 1).Value = 'Beauty'
    $ws1.Cells.Item(7, 2).Value = '2023-06-01'
    $ws1.Cells.Item(7, 3).Value = 25
    $ws1.Cells.Item(7, 4).Value = 10995
    $ws1.Cells.Item(7, 5).Value = 66
    $ws1.Cells.Item(7, 6).Value = 0.0766
    $ws1.Cells.Item(7, 7).Value = 0.0856
    $ws1.Cells.Item(7, 8).Value = 0.2995
    $ws1.Cells.Item(7, 9).Value = 0.335
    $ws1.Cells.Item(8, 1).Value = 'Beauty'
    $ws1.Cells.Item(8, 2).Value = '2023-07-01'
    $ws1.Cells.Item(8, 3).Value = 27
    $ws1.Cells.Item(8, 4).Value = 16090
    $ws1.Cells.Item(8, 5).Value = 70
    $ws1.Cells.Item(8, 6).Value = 0.1121
    $ws1.Cells.Item(8, 7).Value = 0.09080000000000001
    $ws1.Cells.Item(8, 8).Value = 0.4537
    $ws1.Cells.Item(8, 9).Value = 0.3977
    $ws1.Cells.Item(9, 1).Value = 'Beauty'
    $ws1.Cells.Item(9, 2).Value = '2023-08-01'
    $ws1.Cells.Item(9, 3).Value = 24
    $ws1.Cells.Item(9, 4).Value = 9790
    $ws1.Cells.Item(9, 5).Value = 62
    $ws1.Cells.Item(9, 6).Value = 0.0682
    $ws1.Cells.Item(9, 7).Value = 0.0804
    $ws1.Cells.Item(9, 8).Value = 0.2649
    $ws1.Cells.Item(9, 9).Value = 0.2731
    $ws1.Cells.Item(10, 1).Value = 'Beauty'
    $ws1.Cells.Item(10, 2).Value = '2023-09-01'
    $ws1.Cells.Item(10, 3).Value = 20
    $ws1.Cells.Item(10, 4).Value = 6320
    $ws1.Cells.Item(10, 5).Value = 50
    $ws1.Cells.Item(10, 6).Value = 0.044
    $ws1.Cells.Item(10, 7).Value = 0.0649
    $ws1.Cells.Item(10, 8).Value = 0.2676
    $ws1.Cells.Item(10, 9).Value = 0.2941
    $ws1.Cells.Item(11, 1).Value = 'Beauty'
    $ws1.Cells.Item(11, 2).Value = '2023-10-01'
    $ws1.Cells.Item(11, 3).Value = 31
    $ws1.Cells.Item(11, 4).Value = 15355
    $ws1.Cells.Item(11, 5).Value = 83
    $ws1.Cells.Item(11, 6).Value = 0.107
    $ws1.Cells.Item(11, 7).Value = 0.1077
    $ws1.Cells.Item(11, 8).Value = 0.3296
    $ws1.Cells.Item(11, 9).Value = 0.3294
    $ws1.Cells.Item(12, 1).Value = 'Beauty'
    $ws1.Cells.Item(12, 2).Value = '2023-11-01'
    $ws1.Cells.Item(12, 3).Value = 25
    $ws1.Cells.Item(12, 4).Value = 9700
    $ws1.Cells.Item(12, 5).Value = 63
    $ws1.Cells.Item(12, 6).Value = 0.06759999999999999
    $ws1.Cells.Item(12, 7).Value = 0.08169999999999999
    $ws1.Cells.Item(12, 8).Value = 0.2778
    $ws1.Cells.Item(12, 9).Value = 0.3073
    $ws1.Cells.Item(13, 1).Value = 'Beauty'
    $ws1.Cells.Item(13, 2).Value = '2023-12-01'
    $ws1.Cells.Item(13, 3).Value = 25
    $ws1.Cells.Item(13, 4).Value = 12400
    $ws1.Cells.Item(13, 5).Value = 62
    $ws1.Cells.Item(13, 6).Value = 0.0864
    $ws1.Cells.Item(13, 7).Value = 0.0804
    $ws1.Cells.Item(13, 8).Value = 0.2775
    $ws1.Cells.Item(13, 9).Value = 0.2995
    $ws1.Cells.Item(14, 1).Value = 'Beauty'
    $ws1.Cells.Item(14, 2).Value = '2024-01-01'
    $ws1.Cells.Item(14, 3).Value = 1
    $ws1.Cells.Item(14, 4).Value = 1500
    $ws1.Cells.Item(14, 5).Value = 3
    $ws1.Cells.Item(14, 6).Value = 0.0105
    $ws1.Cells.Item(14, 7).Value = 0.0039
    $ws1.Cells.Item(14, 8).Value = 0.9804
    $ws1.Cells.Item(14, 9).Value = 0.75
    $ws1.Cells.Item(15, 1).Value = 'Clothing'
    $ws1.Cells.Item(15, 2).Value = '2023-01-01'
    $ws1.Cells.Item(15, 3).Value = 26
    $ws1.Cells.Item(15, 4).Value = 13125
    $ws1.Cells.Item(15, 5).Value = 72
    $ws1.Cells.Item(15, 6).Value = 0.0844
    $ws1.Cells.Item(15, 7).Value = 0.0805
    $ws1.Cells.Item(15, 8).Value = 0.3702
    $ws1.Cells.Item(15, 9).Value = 0.3692
    $ws1.Cells.Item(16, 1).Value = 'Clothing'
    $ws1.Cells.Item(16, 2).Value = '2023-02-01'
    $ws1.Cells.Item(16, 3).Value = 33
    $ws1.Cells.Item(16, 4).Value = 14560
    $ws1.Cells.Item(16, 5).Value = 75
    $ws1.Cells.Item(16, 6).Value = 0.0936
    $ws1.Cells.Item(16, 7).Value = 0.0839
    $ws1.Cells.Item(16, 8).Value = 0.3305
    $ws1.Cells.Item(16, 9).Value = 0.3505
    $ws1.Cells.Item(17, 1).Value = 'Clothing'
    $ws1.Cells.Item(17, 2).Value = '2023-03-01'
    $ws1.Cells.Item(17, 3).Value = 38
    $ws1.Cells.Item(17, 4).Value = 15065
    $ws1.Cells.Item(17, 5).Value = 111
    $ws1.Cells.Item(17, 6).Value = 0.0968
    $ws1.Cells.Item(17, 7).Value = 0.1242
    $ws1.Cells.Item(17, 8).Value = 0.5197000000000001
    $ws1.Cells.Item(17, 9).Value = 0.5722
    $ws1.Cells.Item(18, 1).Value = 'Clothing'
    $ws1.Cells.Item(18, 2).Value = '2023-04-01'
    $ws1.Cells.Item(18, 3).Value = 36
    $ws1.Cells.Item(18, 4).Value = 13940
    $ws1.Cells.Item(18, 5).Value = 93
    $ws1.Cells.Item(18, 6).Value = 0.0896
    $ws1.Cells.Item(18, 7).Value = 0.104
    $ws1.Cells.Item(18, 8).Value = 0.4116
    $ws1.Cells.Item(18, 9).Value = 0.4346
    $ws1.Cells.Item(19, 1).Value = 'Clothing'
    $ws1.Cells.Item(19, 2).Value = '2023-05-01'
    $ws1.Cells.Item(19, 3).Value = 37
    $ws1.Cells.Item(19, 4).Value = 17455
    $ws1.Cells.Item(19, 5).Value = 97
    $ws1.Cells.Item(19, 6).Value = 0.1122
    $ws1.Cells.Item(19, 7).Value = 0.1085
    $ws1.Cells.Item(19, 8).Value = 0.3284
    $ws1.Cells.Item(19, 9).Value = 0.3745
    $ws1.Cells.Item(20, 1).Value = 'Clothing'
    $ws1.Cells.Item(20, 2).Value = '2023-06-01'
    $ws1.Cells.Item(20, 3).Value = 28
    $ws1.Cells.Item(20, 4).Value = 10170
    $ws1.Cells.Item(20, 5).Value = 67
    $ws1.Cells.Item(20, 6).Value = 0.0654
    $ws1.Cells.Item(20, 7).Value = 0.07489999999999999
    $ws1.Cells.Item(20, 8).Value = 0.277
    $ws1.Cells.Item(20, 9).Value = 0.3401
    $ws1.Cells.Item(21, 1).Value = 'Clothing'
    $ws1.Cells.Item(21, 2).Value = '2023-07-01'
    $ws1.Cells.Item(21, 3).Value = 19
    $ws1.Cells.Item(21, 4).Value = 8250
    $ws1.Cells.Item(21, 5).Value = 45
    $ws1.Cells.Item(21, 6).Value = 0.053
    $ws1.Cells.Item(21, 7).Value = 0.0503
    $ws1.Cells.Item(21, 8).Value = 0.2326
    $ws1.Cells.Item(21, 9).Value = 0.2557
    $ws1.Cells.Item(22, 1).Value = 'Clothing'
    $ws1.Cells.Item(22, 2).Value = '2023-08-01'
    $ws1.Cells.Item(22, 3).Value = 32
    $ws1.Cells.Item(22, 4).Value = 12455
    $ws1.Cells.Item(22, 5).Value = 78
    $ws1.Cells.Item(22, 6).Value = 0.0801
    $ws1.Cells.Item(22, 7).Value = 0.0872
    $ws1.Cells.Item(22, 8).Value = 0.337
    $ws1.Cells.Item(22, 9).Value = 0.3436
    $ws1.Cells.Item(23, 1).Value = 'Clothing'
    $ws1.Cells.Item(23, 2).Value = '2023-09-01'
    $ws1.Cells.Item(23, 3).Value = 20
    $ws1.Cells.Item(23, 4).Value = 9975
    $ws1.Cells.Item(23, 5).Value = 60
    $ws1.Cells.Item(23, 6).Value = 0.0641
    $ws1.Cells.Item(23, 7).Value = 0.06710000000000001
    $ws1.Cells.Item(23, 8).Value = 0.4223
    $ws1.Cells.Item(23, 9).Value = 0.3529
    $ws1.Cells.Item(24, 1).Value = 'Clothing'
    $ws1.Cells.Item(24, 2).Value = '2023-10-01'
    $ws1.Cells.Item(24, 3).Value = 30
    $ws1.Cells.Item(24, 4).Value = 13315
    $ws1.Cells.Item(24, 5).Value = 74
    $ws1.Cells.Item(24, 6).Value = 0.0856
    $ws1.Cells.Item(24, 7).Value = 0.0828
    $ws1.Cells.Item(24, 8).Value = 0.2859
    $ws1.Cells.Item(24, 9).Value = 0.2937
    $ws1.Cells.Item(25, 1).Value = 'Clothing'
    $ws1.Cells.Item(25, 2).Value = '2023-11-01'
    $ws1.Cells.Item(25, 3).Value = 26
    $ws1.Cells.Item(25, 4).Value = 15200
    $ws1.Cells.Item(25, 5).Value = 69
    $ws1.Cells.Item(25, 6).Value = 0.0977
    $ws1.Cells.Item(25, 7).Value = 0.0772
    $ws1.Cells.Item(25, 8).Value = 0.4353
    $ws1.Cells.Item(25, 9).Value = 0.3366
    $ws1.Cells.Item(26, 1).Value = 'Clothing'
    $ws1.Cells.Item(26, 2).Value = '2023-12-01'
    $ws1.Cells.Item(26, 3).Value = 26
    $ws1.Cells.Item(26, 4).Value = 12070
    $ws1.Cells.Item(26, 5).Value = 53
    $ws1.Cells.Item(26, 6).Value = 0.0776
    $ws1.Cells.Item(26, 7).Value = 0.0593
    $ws1.Cells.Item(26, 8).Value = 0.2701
    $ws1.Cells.Item(26, 9).Value = 0.256
    $ws1.Cells.Item(27, 1).Value = 'Electronics'
    $ws1.Cells.Item(27, 2).Value = '2023-01-01'
    $ws1.Cells.Item(27, 3).Value = 25
    $ws1.Cells.Item(27, 4).Value = 9895
    $ws1.Cells.Item(27, 5).Value = 64
    $ws1.Cells.Item(27, 6).Value = 0.0631
    $ws1.Cells.Item(27, 7).Value = 0.07539999999999999
    $ws1.Cells.Item(27, 8).Value = 0.2791
    $ws1.Cells.Item(27, 9).Value = 0.3282
    $ws1.Cells.Item(28, 1).Value = 'Electronics'
    $ws1.Cells.Item(28, 2).Value = '2023-02-01'
    $ws1.Cells.Item(28, 3).Value = 26
    $ws1.Cells.Item(28, 4).Value = 15465
    $ws1.Cells.Item(28, 5).Value = 71
    $ws1.Cells.Item(28, 6).Value = 0.09859999999999999
    $ws1.Cells.Item(28, 7).Value = 0.08359999999999999
    $ws1.Cells.Item(28, 8).Value = 0.351
    $ws1.Cells.Item(28, 9).Value = 0.3318
    $ws1.Cells.Item(29, 1).Value = 'Electronics'
    $ws1.Cells.Item(29, 2).Value = '2023-03-01'
    $ws1.Cells.Item(29, 3).Value = 14
    $ws1.Cells.Item(29, 4).Value = 3380
    $ws1.Cells.Item(29, 5).Value = 32
    $ws1.Cells.Item(29, 6).Value = 0.0215
    $ws1.Cells.Item(29, 7).Value = 0.0377
    $ws1.Cells.Item(29, 8).Value = 0.1166
    $ws1.Cells.Item(29, 9).Value = 0.1649
    $ws1.Cells.Item(30, 1).Value = 'Electronics'
    $ws1.Cells.Item(30, 2).Value = '2023-04-01'
    $ws1.Cells.Item(30, 3).Value = 21
    $ws1.Cells.Item(30, 4).Value = 8025
    $ws1.Cells.Item(30, 5).Value = 52
    $ws1.Cells.Item(30, 6).Value = 0.0511
    $ws1.Cells.Item(30, 7).Value = 0.0612
    $ws1.Cells.Item(30, 8).Value = 0.2369
    $ws1.Cells.Item(30, 9).Value = 0.243
    $ws1.Cells.Item(31, 1).Value = 'Electronics'
    $ws1.Cells.Item(31, 2).Value = '2023-05-01'
    $ws1.Cells.Item(31, 3).Value = 40
    $ws1.Cells.Item(31, 4).Value = 23245
    $ws1.Cells.Item(31, 5).Value = 97
    $ws1.Cells.Item(31, 6).Value = 0.1481
    $ws1.Cells.Item(31, 7).Value = 0.1143
    $ws1.Cells.Item(31, 8).Value = 0.4373
    $ws1.Cells.Item(31, 9).Value = 0.3745
    $ws1.Cells.Item(32, 1).Value = 'Electronics'
    $ws1.Cells.Item(32, 2).Value = '2023-06-01'
    $ws1.Cells.Item(32, 3).Value = 24
    $ws1.Cells.Item(32, 4).Value = 15550
    $ws1.Cells.Item(32, 5).Value = 64
    $ws1.Cells.Item(32, 6).Value = 0.09909999999999999
    $ws1.Cells.Item(32, 7).Value = 0.07539999999999999
    $ws1.Cells.Item(32, 8).Value = 0.4235
    $ws1.Cells.Item(32, 9).Value = 0.3249
    $ws1.Cells.Item(33, 1).Value = 'Electronics'
    $ws1.Cells.Item(33, 2).Value = '2023-07-01'
    $ws1.Cells.Item(33, 3).Value = 26
    $ws1.Cells.Item(33, 4).Value = 11125
    $ws1.Cells.Item(33, 5).Value = 61
    $ws1.Cells.Item(33, 6).Value = 0.0709
    $ws1.Cells.Item(33, 7).Value = 0.0718
    $ws1.Cells.Item(33, 8).Value = 0.3137
    $ws1.Cells.Item(33, 9).Value = 0.3466
    $ws1.Cells.Item(34, 1).Value = 'Electronics'
    $ws1.Cells.Item(34, 2).Value = '2023-08-01'
    $ws1.Cells.Item(34, 3).Value = 38
    $ws1.Cells.Item(34, 4).Value = 14715
    $ws1.Cells.Item(34, 5).Value = 87
    $ws1.Cells.Item(34, 6).Value = 0.09379999999999999
    $ws1.Cells.Item(34, 7).Value = 0.1025
    $ws1.Cells.Item(34, 8).Value = 0.3981
    $ws1.Cells.Item(34, 9).Value = 0.3833
    $ws1.Cells.Item(35, 1).Value = 'Electronics'
    $ws1.Cells.Item(35, 2).Value = '2023-09-01'
    $ws1.Cells.Item(35, 3).Value = 25
    $ws1.Cells.Item(35, 4).Value = 7325
    $ws1.Cells.Item(35, 5).Value = 60
    $ws1.Cells.Item(35, 6).Value = 0.0467
    $ws1.Cells.Item(35, 7).Value = 0.0707
    $ws1.Cells.Item(35, 8).Value = 0.3101
    $ws1.Cells.Item(35, 9).Value = 0.3529
    $ws1.Cells.Item(36, 1).Value = 'Electronics'
    $ws1.Cells.Item(36, 2).Value = '2023-10-01'
    $ws1.Cells.Item(36, 3).Value = 35
    $ws1.Cells.Item(36, 4).Value = 17910
    $ws1.Cells.Item(36, 5).Value = 95
    $ws1.Cells.Item(36, 6).Value = 0.1141
    $ws1.Cells.Item(36, 7).Value = 0.1119
    $ws1.Cells.Item(36, 8).Value = 0.3845
    $ws1.Cells.Item(36, 9).Value = 0.377
    $ws1.Cells.Item(37, 1).Value = 'Electronics'
    $ws1.Cells.Item(37, 2).Value = '2023-11-01'
    $ws1.Cells.Item(37, 3).Value = 27
    $ws1.Cells.Item(37, 4).Value = 10020
    $ws1.Cells.Item(37, 5).Value = 73
    $ws1.Cells.Item(37, 6).Value = 0.0639
    $ws1.Cells.Item(37, 7).Value = 0.08599999999999999
    $ws1.Cells.Item(37, 8).Value = 0.2869
    $ws1.Cells.Item(37, 9).Value = 0.3561
    $ws1.Cells.Item(38, 1).Value = 'Electronics'
    $ws1.Cells.Item(38, 2).Value = '2023-12-01'
    $ws1.Cells.Item(38, 3).Value = 40
    $ws1.Cells.Item(38, 4).Value = 20220
    $ws1.Cells.Item(38, 5).Value = 92
    $ws1.Cells.Item(38, 6).Value = 0.1289
    $ws1.Cells.Item(38, 7).Value = 0.1084
    $ws1.Cells.Item(38, 8).Value = 0.4525
    $ws1.Cells.Item(38, 9).Value = 0.4444
    $ws1.Cells.Item(39, 1).Value = 'Electronics'
    $ws1.Cells.Item(39, 2).Value = '2024-01-01'
    $ws1.Cells.Item(39, 3).Value = 1
    $ws1.Cells.Item(39, 4).Value = 30
    $ws1.Cells.Item(39, 5).Value = 1
    $ws1.Cells.Item(39, 6).Value = 0.0002
    $ws1.Cells.Item(39, 7).Value = 0.0012
    $ws1.Cells.Item(39, 8).Value = 0.0196
    $ws1.Cells.Item(39, 9).Value = 0.25

# Header row styling: bold font, thin border all around, centered
# horizontally, top-aligned vertically (matches the workbook's existing
# header style used on the other generated-report sheets).
$hdr1 = $ws1.Range("A1:I1")
$hdr1.Font.Bold = $true
$hdr1.HorizontalAlignment = -4108
$hdr1.VerticalAlignment = -4160
$hdr1.Borders.LineStyle = 1

# --- New sheet: 08_view_sales_and_customers_mom ----------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "08_view_sales_and_customers_mom"

# Column A (month_start) holds literal text dates -> force Text format.
$ws2.Range("A1:A14").NumberFormat = "@"

    $ws2.Cells.Item(1, 1).Value = 'month_start'
    $ws2.Cells.Item(1, 2).Value = 'year'
    $ws2.Cells.Item(1, 3).Value = 'month'
    $ws2.Cells.Item(1, 4).Value = 'total_revenue'
    $ws2.Cells.Item(1, 5).Value = 'total_units'
    $ws2.Cells.Item(1, 6).Value = 'unique_customers'
    $ws2.Cells.Item(1, 7).Value = 'mom_revenue_growth_pct'
    $ws2.Cells.Item(1, 8).Value = 'ytd_revenue'
    $ws2.Cells.Item(1, 9).Value = 'ytd_units'
    $ws2.Cells.Item(2, 1).Value = '2023-01-01'
    $ws2.Cells.Item(2, 2).Value = 2023
    $ws2.Cells.Item(2, 3).Value = 1
    $ws2.Cells.Item(2, 4).Value = 35450
    $ws2.Cells.Item(2, 5).Value = 195
    $ws2.Cells.Item(2, 6).Value = 76
    $ws2.Cells.Item(2, 8).Value = 35450
    $ws2.Cells.Item(2, 9).Value = 195
    $ws2.Cells.Item(3, 1).Value = '2023-02-01'
    $ws2.Cells.Item(3, 2).Value = 2023
    $ws2.Cells.Item(3, 3).Value = 2
    $ws2.Cells.Item(3, 4).Value = 44060
    $ws2.Cells.Item(3, 5).Value = 214
    $ws2.Cells.Item(3, 6).Value = 85
    $ws2.Cells.Item(3, 7).Value = 24.29
    $ws2.Cells.Item(3, 8).Value = 79510
    $ws2.Cells.Item(3, 9).Value = 409
    $ws2.Cells.Item(4, 1).Value = '2023-03-01'
    $ws2.Cells.Item(4, 2).Value = 2023
    $ws2.Cells.Item(4, 3).Value = 3
    $ws2.Cells.Item(4, 4).Value = 28990
    $ws2.Cells.Item(4, 5).Value = 194
    $ws2.Cells.Item(4, 6).Value = 73
    $ws2.Cells.Item(4, 7).Value = -34.2
    $ws2.Cells.Item(4, 8).Value = 108500
    $ws2.Cells.Item(4, 9).Value = 603
    $ws2.Cells.Item(5, 1).Value = '2023-04-01'
    $ws2.Cells.Item(5, 2).Value = 2023
    $ws2.Cells.Item(5, 3).Value = 4
    $ws2.Cells.Item(5, 4).Value = 33870
    $ws2.Cells.Item(5, 5).Value = 214
    $ws2.Cells.Item(5, 6).Value = 86
    $ws2.Cells.Item(5, 7).Value = 16.83
    $ws2.Cells.Item(5, 8).Value = 142370
    $ws2.Cells.Item(5, 9).Value = 817
    $ws2.Cells.Item(6, 1).Value = '2023-05-01'
    $ws2.Cells.Item(6, 2).Value = 2023
    $ws2.Cells.Item(6, 3).Value = 5
    $ws2.Cells.Item(6, 4).Value = 53150
    $ws2.Cells.Item(6, 5).Value = 259
    $ws2.Cells.Item(6, 6).Value = 105
    $ws2.Cells.Item(6, 7).Value = 56.92
    $ws2.Cells.Item(6, 8).Value = 195520
    $ws2.Cells.Item(6, 9).Value = 1076
    $ws2.Cells.Item(7, 1).Value = '2023-06-01'
    $ws2.Cells.Item(7, 2).Value = 2023
    $ws2.Cells.Item(7, 3).Value = 6
    $ws2.Cells.Item(7, 4).Value = 36715
    $ws2.Cells.Item(7, 5).Value = 197
    $ws2.Cells.Item(7, 6).Value = 77
    $ws2.Cells.Item(7, 7).Value = -30.92
    $ws2.Cells.Item(7, 8).Value = 232235
    $ws2.Cells.Item(7, 9).Value = 1273
    $ws2.Cells.Item(8, 1).Value = '2023-07-01'
    $ws2.Cells.Item(8, 2).Value = 2023
    $ws2.Cells.Item(8, 3).Value = 7
    $ws2.Cells.Item(8, 4).Value = 35465
    $ws2.Cells.Item(8, 5).Value = 176
    $ws2.Cells.Item(8, 6).Value = 72
    $ws2.Cells.Item(8, 7).Value = -3.4
    $ws2.Cells.Item(8, 8).Value = 267700
    $ws2.Cells.Item(8, 9).Value = 1449
    $ws2.Cells.Item(9, 1).Value = '2023-08-01'
    $ws2.Cells.Item(9, 2).Value = 2023
    $ws2.Cells.Item(9, 3).Value = 8
    $ws2.Cells.Item(9, 4).Value = 36960
    $ws2.Cells.Item(9, 5).Value = 227
    $ws2.Cells.Item(9, 6).Value = 94
    $ws2.Cells.Item(9, 7).Value = 4.22
    $ws2.Cells.Item(9, 8).Value = 304660
    $ws2.Cells.Item(9, 9).Value = 1676
    $ws2.Cells.Item(10, 1).Value = '2023-09-01'
    $ws2.Cells.Item(10, 2).Value = 2023
    $ws2.Cells.Item(10, 3).Value = 9
    $ws2.Cells.Item(10, 4).Value = 23620
    $ws2.Cells.Item(10, 5).Value = 170
    $ws2.Cells.Item(10, 6).Value = 65
    $ws2.Cells.Item(10, 7).Value = -36.09
    $ws2.Cells.Item(10, 8).Value = 328280
    $ws2.Cells.Item(10, 9).Value = 1846
    $ws2.Cells.Item(11, 1).Value = '2023-10-01'
    $ws2.Cells.Item(11, 2).Value = 2023
    $ws2.Cells.Item(11, 3).Value = 10
    $ws2.Cells.Item(11, 4).Value = 46580
    $ws2.Cells.Item(11, 5).Value = 252
    $ws2.Cells.Item(11, 6).Value = 96
    $ws2.Cells.Item(11, 7).Value = 97.20999999999999
    $ws2.Cells.Item(11, 8).Value = 374860
    $ws2.Cells.Item(11, 9).Value = 2098
    $ws2.Cells.Item(12, 1).Value = '2023-11-01'
    $ws2.Cells.Item(12, 2).Value = 2023
    $ws2.Cells.Item(12, 3).Value = 11
    $ws2.Cells.Item(12, 4).Value = 34920
    $ws2.Cells.Item(12, 5).Value = 205
    $ws2.Cells.Item(12, 6).Value = 78
    $ws2.Cells.Item(12, 7).Value = -25.03
    $ws2.Cells.Item(12, 8).Value = 409780
    $ws2.Cells.Item(12, 9).Value = 2303
    $ws2.Cells.Item(13, 1).Value = '2023-12-01'
    $ws2.Cells.Item(13, 2).Value = 2023
    $ws2.Cells.Item(13, 3).Value = 12
    $ws2.Cells.Item(13, 4).Value = 44690
    $ws2.Cells.Item(13, 5).Value = 207
    $ws2.Cells.Item(13, 6).Value = 91
    $ws2.Cells.Item(13, 7).Value = 27.98
    $ws2.Cells.Item(13, 8).Value = 454470
    $ws2.Cells.Item(13, 9).Value = 2510
    $ws2.Cells.Item(14, 1).Value = '2024-01-01'
    $ws2.Cells.Item(14, 2).Value = 2024
    $ws2.Cells.Item(14, 3).Value = 1
    $ws2.Cells.Item(14, 4).Value = 1530
    $ws2.Cells.Item(14, 5).Value = 4
    $ws2.Cells.Item(14, 6).Value = 2
    $ws2.Cells.Item(14, 7).Value = -96.58
    $ws2.Cells.Item(14, 8).Value = 1530
    $ws2.Cells.Item(14, 9).Value = 4

$hdr2 = $ws2.Range("A1:I1")
$hdr2.Font.Bold = $true
$hdr2.HorizontalAlignment = -4108
$hdr2.VerticalAlignment = -4160
$hdr2.Borders.LineStyle = 1

Write-Output "07/08 sheets written: $($wb.Worksheets.Count) total worksheets"
